$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "42.888.94"
Set-TextCell "E2" "  +0.20%  "

Set-TextCell "D3" "2.534.54"
Set-TextCell "E3" "  -1.19%  "

Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.15%  "

Set-TextCell "D5" "312.07"
Set-TextCell "E5" "  +0.52%  "

Set-TextCell "D6" "101.35"
Set-TextCell "E6" "  +3.13%  "

Set-TextCell "E7" "  -0.75%  "

Set-TextCell "E8" "  -0.05%  "

Set-TextCell "D9" "0.524"
Set-TextCell "E9" "  -1.10%  "

Set-TextCell "D10" "35.93"
Set-TextCell "E10" "  +0.84%  "

Set-TextCell "D11" "0.0807"
Set-TextCell "E11" "  +0.13%  "

Set-TextCell "D12" "7.37"
Set-TextCell "E12" "  -0.89%  "

Set-TextCell "E13" "  +1.31%  "

Set-TextCell "D14" "2.921.02"
Set-TextCell "E14" "  -1.22%  "

Set-TextCell "D15" "15.42"
Set-TextCell "E15" "  -2.67%  "

Set-TextCell "D16" "2.555.12"
Set-TextCell "E16" "  -1.91%  "

Set-TextCell "D17" "0.819"
Set-TextCell "E17" "  -2.28%  "

Set-TextCell "D18" "42.844.26"
Set-TextCell "E18" "  +0.03%  "

Set-TextCell "D19" "6.68"
Set-TextCell "E19" "  -0.61%  "

Set-TextCell "D20" "12.44"
Set-TextCell "E20" "  +0.75%  "

Set-TextCell "D21" "0.0₃0955"
Set-TextCell "E21" "  -0.43%  "

Set-TextCell "E22" "  +0.90%  "

Set-TextCell "D23" "244.34"
Set-TextCell "E23" "  -1.31%  "

Set-TextCell "D24" "2.89"
Set-TextCell "E24" "  -1.17%  "

Set-TextCell "E25" "  -0.37%  "

Set-TextCell "E26" "  +0.15%  "

Set-TextCell "D27" "25.62"
Set-TextCell "E27" "  -4.98%  "

Set-TextCell "E28" "  -2.19%  "

Set-TextCell "D29" "10.23"
Set-TextCell "E29" "  +0.78%  "

Set-TextCell "D30" "38.80"
Set-TextCell "E30" "  -2.62%  "

Set-TextCell "D31" "161.34"
Set-TextCell "E31" "  +1.91%  "

Set-TextCell "D32" "5.84"
Set-TextCell "E32" "  +1.54%  "

Set-TextCell "D33" "2.79"
Set-TextCell "E33" "  +7.86%  "

Set-TextCell "E34" "  +0.12%  "

Set-TextCell "D35" "0.0792"
Set-TextCell "E35" "  -0.43%  "

Set-TextCell "D36" "18.44"
Set-TextCell "E36" "  -1.39%  "

Set-TextCell "D37" "3.12"
Set-TextCell "E37" "  -5.17%  "

Set-TextCell "D38" "1.97"
Set-TextCell "E38" "  -5.51%  "

Set-TextCell "E39" "  -0.21%  "

Set-TextCell "D40" "0.118"
Set-TextCell "E40" "  +0.13%  "

Set-TextCell "D41" "4.19"
Set-TextCell "E41" "  +2.50%  "

Set-TextCell "D42" "22.03"
Set-TextCell "E42" "  -4.17%  "

Set-TextCell "B43" "NEARProtocol"
Set-TextCell "C43" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D43" "3.34"
Set-TextCell "E43" "  +4.60%  "

Set-TextCell "B44" "FirstDigitalUSD"
Set-TextCell "C44" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D44" "1.00"
Set-TextCell "E44" "  +0.07%  "

Set-TextCell "D45" "0.0300"
Set-TextCell "E45" "  -0.39%  "

Set-TextCell "D46" "2.000.19"
Set-TextCell "E46" "  +0.52%  "

Set-TextCell "D47" "9.17"
Set-TextCell "E47" "  +1.95%  "

Set-TextCell "D48" "2.772.82"
Set-TextCell "E48" "  -1.34%  "

Set-TextCell "E49" "  -1.08%  "

Set-TextCell "D50" "79.83"
Set-TextCell "E50" "  -1.60%  "

Set-TextCell "D51" "72.65"
Set-TextCell "E51" "  -1.08%  "
